$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two "Requisitos" entries (rows 23 and 24) need to swap places:
# row 23 should now show the LOT2028 text, row 24 the LOT2052 text.
$lot2028 = "LOT2028 -  Tecnologia de Processos Fermentativos  (Requisito fraco)`n"
$lot2052 = "LOT2052 -  Tecnologia de Bebidas Experimental  (Indicação de Conjunto)`n"

$ws.Range("B23").Value = $lot2028
$ws.Range("C23").Value = $lot2028

$ws.Range("B24").Value = $lot2052
$ws.Range("C24").Value = $lot2052
